$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row before the current last data row (row 22) ---
# This pushes the old row 22 (and the trailing signature rows 27/28) down by
# one row, matching the new layout (data rows 16-23, signature rows 28-29).
$ws.Rows(22).Insert()

# Copy the formatting of the row above (21) into the freshly inserted row
# (22) so the new row matches the style used by the other "middle" data rows.
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header / summary fields ---
$ws.Range("E11").Value = 377867
$ws.Range("F13").Value = 8

# --- Data rows 16-23: Tipo Doc / N Doc / Nombre / Periodo / Valor Mora / Salario Basico ---
$periods = @("2507","2506","2505","2504","2503","2502","2501","2412")
$valorMora = @(52000,52000,52000,52000,52000,52000,52000,13867)
$salarioBasico = @(908526,908526,908526,908526,908526,908526,908526,908526)

for ($i = 0; $i -lt 8; $i++) {
  $r = 16 + $i
  $ws.Cells.Item($r, 2).Value = "CC"
  $ws.Cells.Item($r, 3).Value = "73578753"
  $ws.Cells.Item($r, 4).Value = "DARWIN ALEXANDER MERCADO VEGA"
  $ws.Cells.Item($r, 5).Value = $periods[$i]
  $ws.Cells.Item($r, 6).Value = $valorMora[$i]
  $ws.Cells.Item($r, 7).Value = $salarioBasico[$i]
}

# --- Signature block (now rows 28-29) ---
$ws.Range("B28").Value = "___________________________________"
$ws.Range("H28").Value = "___________________________________"
$ws.Range("B29").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H29").Value = "FIRMA DEL REPRESENTANTE LEGAL"
